# Horarioatencion.docx: the schedule table's day-of-week cell changes
# from "MIERCOLES" (Wednesday) to "LUNES" (Monday).
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "MIERCOLES",  # FindText
    $true,        # MatchCase
    $true,        # MatchWholeWord
    $false,       # MatchWildcards
    $false,       # MatchSoundsLike
    $false,       # MatchAllWordForms
    $true,        # Forward
    1,            # Wrap (wdFindContinue)
    $false,       # Format
    "LUNES",      # ReplaceWith
    2             # Replace (wdReplaceAll)
)
